$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(95, 8).Value = 63992.5  # H95: 63995 -> 63992.5
$ws.Cells.Item(95, 10).Value = 63992.5  # J95: 63995 -> 63992.5
$ws.Cells.Item(95, 12).Value = 63992.5  # L95: 63995 -> 63992.5
$ws.Cells.Item(95, 14).Value = -69484.5  # N95: -69487 -> -69484.5
$ws.Cells.Item(103, 8).Value = 833.3333  # H103: 979.8 -> 833.3333
$ws.Cells.Item(103, 9).Value = 674.5  # I103: 799 -> 674.5
$ws.Cells.Item(103, 10).Value = 1151  # J103: 1100.3334 -> 1151
$ws.Cells.Item(103, 11).Value = 2023.5  # K103: 2397 -> 2023.5
$ws.Cells.Item(103, 12).Value = 3453  # L103: 3301.0002 -> 3453
$ws.Cells.Item(103, 13).Value = -1437.5  # M103: -1811 -> -1437.5
$ws.Cells.Item(103, 14).Value = -4625  # N103: -4473.0002 -> -4625
$ws.Cells.Item(105, 8).Value = 87000  # H105: 0 -> 87000
$ws.Cells.Item(105, 10).Value = 87000  # J105: 0 -> 87000
$ws.Cells.Item(105, 12).Value = 87000  # L105: 0 -> 87000
$ws.Cells.Item(105, 14).Value = -93988  # N105: <MISSING> -> -93988
$ws.Cells.Item(137, 8).Value = 4462.6523  # H137: 4539.9556 -> 4462.6523
$ws.Cells.Item(137, 9).Value = 2211.7778  # I137: 2259 -> 2211.7778
$ws.Cells.Item(137, 11).Value = 6635.3334  # K137: 6777 -> 6635.3334
$ws.Cells.Item(137, 13).Value = -4085.3334  # M137: -4227 -> -4085.3334
$ws.Cells.Item(138, 8).Value = 2720.4  # H138: 2737.31 -> 2720.4
$ws.Cells.Item(138, 10).Value = 2741.6836  # J138: 2758.9387 -> 2741.6836
$ws.Cells.Item(138, 12).Value = 8225.050799999999  # L138: 8276.8161 -> 8225.050799999999
$ws.Cells.Item(138, 14).Value = -18505.0508  # N138: -18556.8161 -> -18505.0508

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 30814.166  # H32: 31513.02 -> 30814.166
$ws.Cells.Item(32, 9).Value = 13782.233  # I32: 13611.517 -> 13782.233
$ws.Cells.Item(32, 10).Value = 59200.723  # J32: 59260.35 -> 59200.723
$ws.Cells.Item(32, 11).Value = 13782.233  # K32: 13611.517 -> 13782.233
$ws.Cells.Item(32, 12).Value = 59200.723  # L32: 59260.35 -> 59200.723
$ws.Cells.Item(32, 13).Value = -13495.233  # M32: -13324.517 -> -13495.233
$ws.Cells.Item(32, 14).Value = -59774.723  # N32: -59834.35 -> -59774.723
$ws.Cells.Item(61, 8).Value = 3774.85  # H61: 3897.205 -> 3774.85
$ws.Cells.Item(61, 9).Value = 2257.1177  # I61: 2358.2 -> 2257.1177
$ws.Cells.Item(61, 10).Value = 4896.6523  # J61: 4859.0835 -> 4896.6523
$ws.Cells.Item(61, 11).Value = 2257.1177  # K61: 2358.2 -> 2257.1177
$ws.Cells.Item(61, 12).Value = 4896.6523  # L61: 4859.0835 -> 4896.6523
$ws.Cells.Item(61, 13).Value = -2045.1177  # M61: -2146.2 -> -2045.1177
$ws.Cells.Item(61, 14).Value = -5320.6523  # N61: -5283.0835 -> -5320.6523
$ws.Cells.Item(80, 8).Value = 99000  # H80: 0 -> 99000
$ws.Cells.Item(80, 10).Value = 99000  # J80: 0 -> 99000
$ws.Cells.Item(80, 12).Value = 99000  # L80: 0 -> 99000
$ws.Cells.Item(80, 14).Value = -100996  # N80: <MISSING> -> -100996
$ws.Cells.Item(83, 8).Value = 99000  # H83: 0 -> 99000
$ws.Cells.Item(83, 10).Value = 99000  # J83: 0 -> 99000
$ws.Cells.Item(83, 12).Value = 297000  # L83: 0 -> 297000
$ws.Cells.Item(83, 14).Value = -306984  # N83: <MISSING> -> -306984
$ws.Cells.Item(102, 8).Value = 3000  # H102: 0 -> 3000
$ws.Cells.Item(102, 9).Value = 3000  # I102: 0 -> 3000
$ws.Cells.Item(102, 11).Value = 3000  # K102: 0 -> 3000
$ws.Cells.Item(102, 13).Value = -1378  # M102: <MISSING> -> -1378
$ws.Cells.Item(136, 8).Value = 3774.85  # H136: 3897.205 -> 3774.85
$ws.Cells.Item(136, 9).Value = 2257.1177  # I136: 2358.2 -> 2257.1177
$ws.Cells.Item(136, 10).Value = 4896.6523  # J136: 4859.0835 -> 4896.6523
$ws.Cells.Item(136, 11).Value = 6771.353099999999  # K136: 7074.599999999999 -> 6771.353099999999
$ws.Cells.Item(136, 12).Value = 14689.9569  # L136: 14577.2505 -> 14689.9569
$ws.Cells.Item(136, 13).Value = -4221.353099999999  # M136: -4524.599999999999 -> -4221.353099999999
$ws.Cells.Item(136, 14).Value = -19789.9569  # N136: -19677.2505 -> -19789.9569

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(54, 8).Value = 6663  # H54: 9494.5 -> 6663
$ws.Cells.Item(54, 10).Value = 1000  # J54: 0 -> 1000
$ws.Cells.Item(54, 12).Value = 1000  # L54: 0 -> 1000
$ws.Cells.Item(54, 14).Value = -1968  # N54: <MISSING> -> -1968
$ws.Cells.Item(88, 8).Value = 26965  # H88: 30167.875 -> 26965
$ws.Cells.Item(88, 10).Value = 27835.625  # J88: 31620.428 -> 27835.625
$ws.Cells.Item(88, 12).Value = 27835.625  # L88: 31620.428 -> 27835.625
$ws.Cells.Item(88, 14).Value = -28647.625  # N88: -32432.428 -> -28647.625
$ws.Cells.Item(91, 8).Value = 26965  # H91: 30167.875 -> 26965
$ws.Cells.Item(91, 10).Value = 27835.625  # J91: 31620.428 -> 27835.625
$ws.Cells.Item(91, 12).Value = 27835.625  # L91: 31620.428 -> 27835.625
$ws.Cells.Item(91, 14).Value = -30643.625  # N91: -34428.428 -> -30643.625
$ws.Cells.Item(99, 8).Value = 2999  # H99: 2999.5 -> 2999
$ws.Cells.Item(99, 9).Value = 2999  # I99: 2999.5 -> 2999
$ws.Cells.Item(99, 11).Value = 2999  # K99: 2999.5 -> 2999
$ws.Cells.Item(99, 13).Value = -1501  # M99: -1501.5 -> -1501
$ws.Cells.Item(105, 8).Value = 4167.5713  # H105: 4325.95 -> 4167.5713
$ws.Cells.Item(105, 9).Value = 3965.125  # I105: 4388.7144 -> 3965.125
$ws.Cells.Item(105, 11).Value = 3965.125  # K105: 4388.7144 -> 3965.125
$ws.Cells.Item(105, 13).Value = -2218.125  # M105: -2641.7144 -> -2218.125

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 10926.7  # H62: 10487.909 -> 10926.7
$ws.Cells.Item(62, 10).Value = 9895  # J62: 8630 -> 9895
$ws.Cells.Item(62, 12).Value = 9895  # L62: 8630 -> 9895
$ws.Cells.Item(62, 14).Value = -11143  # N62: -9878 -> -11143
$ws.Cells.Item(65, 8).Value = 10926.7  # H65: 10487.909 -> 10926.7
$ws.Cells.Item(65, 10).Value = 9895  # J65: 8630 -> 9895
$ws.Cells.Item(65, 12).Value = 49475  # L65: 43150 -> 49475
$ws.Cells.Item(65, 14).Value = -55715  # N65: -49390 -> -55715
$ws.Cells.Item(88, 8).Value = 30199.666  # H88: 30937.125 -> 30199.666
$ws.Cells.Item(88, 10).Value = 30199.666  # J88: 30937.125 -> 30199.666
$ws.Cells.Item(88, 12).Value = 30199.666  # L88: 30937.125 -> 30199.666
$ws.Cells.Item(88, 14).Value = -31011.666  # N88: -31749.125 -> -31011.666
$ws.Cells.Item(91, 8).Value = 30199.666  # H91: 30937.125 -> 30199.666
$ws.Cells.Item(91, 10).Value = 30199.666  # J91: 30937.125 -> 30199.666
$ws.Cells.Item(91, 12).Value = 30199.666  # L91: 30937.125 -> 30199.666
$ws.Cells.Item(91, 14).Value = -33007.666  # N91: -33745.125 -> -33007.666
$ws.Cells.Item(92, 8).Value = 25549.25  # H92: 34300.5 -> 25549.25
$ws.Cells.Item(92, 10).Value = 25549.25  # J92: 34300.5 -> 25549.25
$ws.Cells.Item(92, 12).Value = 25549.25  # L92: 34300.5 -> 25549.25
$ws.Cells.Item(92, 14).Value = -30541.25  # N92: -39292.5 -> -30541.25
$ws.Cells.Item(95, 8).Value = 19997.5  # H95: 20000 -> 19997.5
$ws.Cells.Item(95, 10).Value = 19997.5  # J95: 20000 -> 19997.5
$ws.Cells.Item(95, 12).Value = 19997.5  # L95: 20000 -> 19997.5
$ws.Cells.Item(95, 14).Value = -25489.5  # N95: -25492 -> -25489.5
$ws.Cells.Item(96, 8).Value = 18274.334  # H96: 0 -> 18274.334
$ws.Cells.Item(96, 10).Value = 18274.334  # J96: 0 -> 18274.334
$ws.Cells.Item(96, 12).Value = 18274.334  # L96: 0 -> 18274.334
$ws.Cells.Item(96, 14).Value = -23766.334  # N96: <MISSING> -> -23766.334
$ws.Cells.Item(105, 8).Value = 1939.08  # H105: 1943.08 -> 1939.08
$ws.Cells.Item(105, 10).Value = 2939.6667  # J105: 2950.7778 -> 2939.6667
$ws.Cells.Item(105, 12).Value = 2939.6667  # L105: 2950.7778 -> 2939.6667
$ws.Cells.Item(105, 14).Value = -6433.6667  # N105: -6444.7778 -> -6433.6667

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 3947.5  # H80: 3994.5 -> 3947.5
$ws.Cells.Item(80, 10).Value = 3947.5  # J80: 3994.5 -> 3947.5
$ws.Cells.Item(80, 12).Value = 11842.5  # L80: 11983.5 -> 11842.5
$ws.Cells.Item(80, 14).Value = -13714.5  # N80: -13855.5 -> -13714.5
$ws.Cells.Item(83, 8).Value = 3947.5  # H83: 3994.5 -> 3947.5
$ws.Cells.Item(83, 10).Value = 3947.5  # J83: 3994.5 -> 3947.5
$ws.Cells.Item(83, 12).Value = 35527.5  # L83: 35950.5 -> 35527.5
$ws.Cells.Item(83, 14).Value = -44887.5  # N83: -45310.5 -> -44887.5
$ws.Cells.Item(122, 8).Value = 1480.4615  # H122: 1290.8422 -> 1480.4615
$ws.Cells.Item(122, 9).Value = 1459.4  # I122: 1143.3636 -> 1459.4
$ws.Cells.Item(122, 11).Value = 13134.6  # K122: 10290.2724 -> 13134.6
$ws.Cells.Item(122, 13).Value = -10684.6  # M122: -7840.2724 -> -10684.6

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 9000  # H5: 10000 -> 9000
$ws.Cells.Item(5, 9).Value = 9000  # I5: 10000 -> 9000
$ws.Cells.Item(5, 11).Value = 9000  # K5: 10000 -> 9000
$ws.Cells.Item(5, 13).Value = -8888  # M5: -9888 -> -8888
$ws.Cells.Item(9, 8).Value = 931.6667  # H9: 503.5 -> 931.6667
$ws.Cells.Item(9, 9).Value = 0  # I9: 7 -> 0
$ws.Cells.Item(9, 10).Value = 931.6667  # J9: 1000 -> 931.6667
$ws.Cells.Item(9, 11).Value = 0  # K9: 7 -> 0
$ws.Cells.Item(9, 12).Value = 931.6667  # L9: 1000 -> 931.6667
$ws.Cells.Item(9, 13).ClearContents()  # M9: was 163
$ws.Cells.Item(9, 14).Value = -1271.6667  # N9: -1340 -> -1271.6667
$ws.Cells.Item(17, 8).Value = 0  # H17: 1970 -> 0
$ws.Cells.Item(17, 9).Value = 0  # I17: 2001 -> 0
$ws.Cells.Item(17, 10).Value = 0  # J17: 1954.5 -> 0
$ws.Cells.Item(17, 11).Value = 0  # K17: 2001 -> 0
$ws.Cells.Item(17, 12).Value = 0  # L17: 1954.5 -> 0
$ws.Cells.Item(17, 13).ClearContents()  # M17: was -1833
$ws.Cells.Item(17, 14).ClearContents()  # N17: was -2290.5
$ws.Cells.Item(122, 8).Value = 3453.6316  # H122: 3462.1667 -> 3453.6316
$ws.Cells.Item(122, 10).Value = 3439.5334  # J122: 3440 -> 3439.5334
$ws.Cells.Item(122, 12).Value = 10318.6002  # L122: 10320 -> 10318.6002
$ws.Cells.Item(122, 14).Value = -15218.6002  # N122: -15220 -> -15218.6002
$ws.Cells.Item(126, 8).Value = 4097.8335  # H126: 4068.2163 -> 4097.8335
$ws.Cells.Item(126, 9).Value = 3467.9092  # I126: 3429.0833 -> 3467.9092
$ws.Cells.Item(126, 11).Value = 10403.7276  # K126: 10287.2499 -> 10403.7276
$ws.Cells.Item(126, 13).Value = -7933.7276  # M126: -7817.249899999999 -> -7933.7276

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 3576.5  # H100: 3645.1333 -> 3576.5
$ws.Cells.Item(100, 9).Value = 3504  # I100: 3647.7144 -> 3504
$ws.Cells.Item(100, 10).Value = 3649  # J100: 3642.875 -> 3649
$ws.Cells.Item(100, 11).Value = 3504  # K100: 3647.7144 -> 3504
$ws.Cells.Item(100, 12).Value = 3649  # L100: 3642.875 -> 3649
$ws.Cells.Item(100, 13).Value = -2963  # M100: -3106.7144 -> -2963
$ws.Cells.Item(100, 14).Value = -4731  # N100: -4724.875 -> -4731
$ws.Cells.Item(132, 8).Value = 3781.6365  # H132: 3929.238 -> 3781.6365
$ws.Cells.Item(132, 9).Value = 2895.25  # I132: 3500 -> 2895.25
$ws.Cells.Item(132, 10).Value = 3978.611  # J132: 3950.7 -> 3978.611
$ws.Cells.Item(132, 11).Value = 8685.75  # K132: 10500 -> 8685.75
$ws.Cells.Item(132, 12).Value = 11935.833  # L132: 11852.1 -> 11935.833
$ws.Cells.Item(132, 13).Value = -6155.75  # M132: -7970 -> -6155.75
$ws.Cells.Item(132, 14).Value = -16995.833  # N132: -16912.1 -> -16995.833

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 13982.571  # H41: 13484.125 -> 13982.571
$ws.Cells.Item(41, 10).Value = 9994.5  # J41: 9994.6 -> 9994.5
$ws.Cells.Item(41, 12).Value = 9994.5  # L41: 9994.6 -> 9994.5
$ws.Cells.Item(41, 14).Value = -10774.5  # N41: -10774.6 -> -10774.5
$ws.Cells.Item(136, 8).Value = 1930.3462  # H136: 2160.3667 -> 1930.3462
$ws.Cells.Item(136, 9).Value = 2030.3125  # I136: 2140.1333 -> 2030.3125
$ws.Cells.Item(136, 10).Value = 1770.4  # J136: 2180.6 -> 1770.4
$ws.Cells.Item(136, 11).Value = 6090.9375  # K136: 6420.3999 -> 6090.9375
$ws.Cells.Item(136, 12).Value = 5311.200000000001  # L136: 6541.799999999999 -> 5311.200000000001
$ws.Cells.Item(136, 13).Value = -3540.9375  # M136: -3870.3999 -> -3540.9375
$ws.Cells.Item(136, 14).Value = -10411.2  # N136: -11641.8 -> -10411.2
